# Add two new parameter columns ("xlabel", "ylabel") into Sheet1, inserted
# immediately to the left of what is currently column T ("auto_deconv").
# Everything at/after column T shifts two places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at T:U (existing T:U and beyond shift to V:W, etc.)
$ws.Range("T1:U1").EntireColumn.Insert()

# New header labels for the inserted columns (row 1)
$ws.Range("T1").Value = "xlabel"
$ws.Range("U1").Value = "ylabel"

# New default values for the inserted columns (row 2)
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 1

# The inserted columns should keep the same width as column S (11 chars)
$ws.Range("T1:U1").ColumnWidth = $ws.Range("S1").ColumnWidth

# Update the active selection to match the saved state of the workbook
$ws.Range("U3").Select()
